# Kupa yolu maçları girildi ve Hubuş FK vs Armedospor maçı kaydı eklendi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Çeyrek Final (quarterfinal) matchups: fill in the real team names ---
$ws.Range("A3").Value = "61.Alay"
$ws.Range("B3").Value = "Ajans Of"

$ws.Range("A6").Value = "Fortuna United"
$ws.Range("B6").Value = "Kural Kesiciler"

$ws.Range("A9").Value = "Çirihtalar"
$ws.Range("B9").Value = "Araklı 1961 Spor"

$ws.Range("A12").Value = "Hubuş FK"
$ws.Range("B12").Value = "Of FK"

# --- Yarı Final / Final / 3.lük rows: clear the leftover placeholder team names ---
$ws.Range("A15:B15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()

$ws.Range("A18:B18").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("B20").ClearContents()

$ws.Range("A21:B21").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("B23").ClearContents()

$ws.Range("A24:B24").ClearContents()

# --- View state: scroll back to top-left and move the active selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D22").Select()
